# Add the "Distributor_Add" sheet (after Sheet1) with distributor data,
# a mailto hyperlink on the email cell, and matching column widths.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Distributor_Add"

# Header row
$ws.Range("A1").Value = "Distributor Code"
$ws.Range("B1").Value = "Distributor Name"
$ws.Range("C1").Value = "Contact person name"
$ws.Range("D1").Value = "Mobile Number"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Address1"
$ws.Range("G1").Value = "Address2"
$ws.Range("H1").Value = "Address3"
$ws.Range("I1").Value = "PinCode"
$ws.Range("J1").Value = "Pan Number"
$ws.Range("K1").Value = "FSSAI Licence Number"

# Data row
$ws.Range("A2").Value = "DB5410"
$ws.Range("B2").Value = "Aditya"
$ws.Range("C2").Value = "Aditya Thakur"
$ws.Range("D2").Value = 9654862012
$ws.Range("E2").Value = "aniket.jadhav@heerasoftware.com"
$ws.Range("F2").Value = "ramai nagar"
$ws.Range("G2").Value = "washing center"
$ws.Range("H2").Value = "narhe"
$ws.Range("I2").Value = 441611
$ws.Range("J2").Value = "POYTF5487K"
$ws.Range("K2").Value = 10012022001234

# Hyperlink on the email cell (applies the built-in Hyperlink style too)
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:aniket.jadhav@heerasoftware.com") | Out-Null

# Column widths (values chosen so the stored xlsx width matches the target)
$ws.Columns.Item(1).ColumnWidth = 13.333333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws.Columns.Item(3).ColumnWidth = 20.166666666666668
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 28.666666666666668
$ws.Columns.Item(6).ColumnWidth = 10.166666666666666
$ws.Columns.Item(7).ColumnWidth = 12.5
$ws.Columns.Item(10).ColumnWidth = 12.333333333333334
$ws.Columns.Item(11).ColumnWidth = 19.833333333333332

# Make the new sheet the active / selected one, with the same selection as
# the source workbook.
$ws.Activate()
$ws.Range("K4").Select()
